# Update trad and label database (encore)
# Retranslate / update the English labels for 8 codes in column B.
# Codes (column A) stay the same and the sheet stays sorted alphabetically
# by code; only the label text in column B changes for these rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Cost of labor"                                   # C_L
$ws.Range("B7").Value  = "Non-qualified employment (in thousands)"         # F_L_NQ
$ws.Range("B8").Value  = "Qualified employment (in thousands)"             # F_L_Q
$ws.Range("B9").Value  = "Employment in thousands (commercial sector)"     # F_L_SM
$ws.Range("B12").Value = "Investment (SNF + EI)"                           # I_SNF
$ws.Range("B16").Value = "Intermediate consumption price"                  # PCI
$ws.Range("B21").Value = "Production price (commercial sector)"            # PYM
$ws.Range("B28").Value = "VA (commercial sector)"                          # VA_SM

# Move the active selection from D10 to B2, matching the refreshed view.
$null = $ws.Range("B2").Select()
